$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing sheet "Hoja1" -> "ExistingHoja1"
$ws.Name = "ExistingHoja1"

# Move/update the selection to B20 (single cell)
$ws.Range("B20").Select()
